$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text in A1 ("UBSG:SWX Quotes" -> "UBS Quotes")
$ws.Range("A1").Value = "UBS Quotes"

# Touch A4 so it materializes as a blank cell (extends used range/dimension to A1:A4)
$ws.Range("A4").Style = "Normal"

# Move/save the active selection to C3
$ws.Range("C3").Select()
